$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F8").Value = 4582
$ws.Range("F11").Value = 1753
$ws.Range("F13").Value = 694
$ws.Range("F18").Value = 1550
$ws.Range("F20").Value = 705
$ws.Range("F34").Value = 4152

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 12
$ws.Range("F27").Value = 241
$ws.Range("F39").Value = 26

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1714
$ws.Range("F8").Value = 185

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1714
$ws.Range("F12").Value = 185
$ws.Range("F15").Value = 12
$ws.Range("F16").Value = 4582
$ws.Range("F19").Value = 1753
$ws.Range("F21").Value = 694
$ws.Range("F27").Value = 1550
$ws.Range("F31").Value = 705
$ws.Range("F44").Value = 241
$ws.Range("F48").Value = 4152
$ws.Range("F50").Value = 26

$wb.Save()
